$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Rename the main sheet and add the two new lookup sheets
# ---------------------------------------------------------------------------
$mainWs = $wb.Worksheets.Item(1)
$mainWs.Name = "-23&P"

$chaptersWs = $wb.Worksheets.Add($null, $mainWs)
$chaptersWs.Name = "Chapters"

$proceduresWs = $wb.Worksheets.Add($null, $chaptersWs)
$proceduresWs.Name = "Procedures"

# ---------------------------------------------------------------------------
# 2. Populate the Chapters sheet + turn it into a table
# ---------------------------------------------------------------------------
$chapters = @(
    "Chapter",
    "Chapter 1",
    "Operator Procedures",
    "Troubleshooting Master Index",
    "Maintainer Troubeshooting",
    "Maintainer PMCS",
    "Maintainer Maintenance",
    "Destruction of Equipment to Prevent Enemy Use",
    "Repair Parts and Special Tools List",
    "Supporting Information",
    "Rear Matter"
)
for ($i = 0; $i -lt $chapters.Length; $i++) {
    $chaptersWs.Cells.Item($i + 1, 1).Value = $chapters[$i]
}
$chaptersWs.Columns.Item(1).ColumnWidth = 43.6328125

$chaptersRange = $chaptersWs.Range("A1:A11")
$chaptersTable = $chaptersWs.ListObjects.Add(1, $chaptersRange, $null, 1)
$chaptersTable.Name = "Table4"
$chaptersTable.TableStyle = "TableStyleMedium9"

# ---------------------------------------------------------------------------
# 3. Populate the Procedures sheet + turn it into a table
# ---------------------------------------------------------------------------
$procedures = @(
    "Procedures",
    "inspect",
    "test",
    "service",
    "adjust",
    "align",
    "calibration",
    "remove",
    "install",
    "replace",
    "repair",
    "paint",
    "overhaul",
    "rebuild",
    "lube",
    "mark",
    "pack",
    "unpack",
    "preservation",
    "prepforuse",
    "assem",
    "disassem",
    "clean",
    "ndi",
    "ris",
    "pis",
    "tow",
    "jack",
    "park",
    "moor",
    "cover",
    "hoist",
    "sling",
    "extpwr",
    "prepstore",
    "prepship",
    "transport",
    "arm",
    "load",
    "unload",
    "installperdev",
    "uninstallperdev",
    "upgrade",
    "configure",
    "debug"
)
for ($i = 0; $i -lt $procedures.Length; $i++) {
    $proceduresWs.Cells.Item($i + 1, 1).Value = $procedures[$i]
}
$proceduresWs.Columns.Item(1).ColumnWidth = 15

$proceduresRange = $proceduresWs.Range("A1:A45")
$proceduresTable = $proceduresWs.ListObjects.Add(1, $proceduresRange, $null, 1)
$proceduresTable.Name = "Table5"
$proceduresTable.TableStyle = "TableStyleMedium9"

# ---------------------------------------------------------------------------
# 4. Shorten the chapter header labels on the main sheet (drop "Chapter N - ")
# ---------------------------------------------------------------------------
$mainWs.Range("C11").Value = "Operator Procedures"
$mainWs.Range("C17").Value = "Troubleshooting Master Index"
$mainWs.Range("C19").Value = "Maintainer Troubeshooting"
$mainWs.Range("C25").Value = "Maintainer PMCS"
$mainWs.Range("C31").Value = "Maintainer Maintenance"
$mainWs.Range("C39").Value = "Destruction of Equipment to Prevent Enemy Use"
$mainWs.Range("C43").Value = "Repair Parts and Special Tools List"
$mainWs.Range("C51").Value = "Supporting Information"

# ---------------------------------------------------------------------------
# 5. Lower-case the procedure keywords in column D (rows 33-38)
# ---------------------------------------------------------------------------
$mainWs.Range("D33").Value = "service"
$mainWs.Range("D34").Value = "install"
$mainWs.Range("D35").Value = "repair"
$mainWs.Range("D36").Value = "prepstore"
$mainWs.Range("D37").Value = "transport"
$mainWs.Range("D38").Value = "clean"

# ---------------------------------------------------------------------------
# 6. Add the data validation drop-downs driven by the new lookup tables
# ---------------------------------------------------------------------------
$procRange = $mainWs.Range("D33:D38")
$procRange.Validation.Delete()
$procRange.Validation.Add(3, 1, 1, "=Procedures!`$A`$2:`$A`$45")
$procRange.Validation.IgnoreBlank = 1
$procRange.Validation.InCellDropdown = 1

$chapterRange = $mainWs.Range("C7,C11,C17,C19,C25,C31,C39,C43,C51,C63")
$chapterRange.Validation.Delete()
$chapterRange.Validation.Add(3, 1, 1, "=Chapters!`$A`$2:`$A`$11")
$chapterRange.Validation.IgnoreBlank = 1
$chapterRange.Validation.InCellDropdown = 1

# ---------------------------------------------------------------------------
# 7. Restore original active sheet / selection
# ---------------------------------------------------------------------------
$mainWs.Activate()
$mainWs.Range("I20").Select()
